$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is forced to text via a leading apostrophe (quote-prefix),
# matching the source workbook where these Price/Volume columns are stored
# as plain text, not numbers (e.g. "34.587.55", "  +1.07%  ").
# Style is reset to Normal afterwards so no stray text-format style sticks
# to the cell (keeps formatting identical to the original cells).

$ws.Range("D2").Value = "'34.587.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.07%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.797.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.76%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'227.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.44%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +1.58%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +1.92%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.91%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0696"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.28%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.25%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.059.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.84%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'11.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.79%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.795.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.05%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.638"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.18%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'34.575.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.18%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +2.87%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'68.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.32%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.0₃0804"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.63%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'247.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.05%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.50%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.12%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +2.29%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'168.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.43%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.37%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.24%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'16.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +2.10%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.00%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +9.62%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +2.24%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.0526"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.96%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.33%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +2.30%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.431.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.99%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +7.62%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +2.78%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.92%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.22%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'84.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +5.49%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.71%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.937"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.14%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +3.32%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'13.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.99%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +3.21%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.47%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.74%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.959.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.81%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'105.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.04%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.03%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -4.09%  "
$ws.Range("E51").Style = "Normal"
